$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.132.74'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.652.02'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E4').Value = '  +0.36%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.24'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5199'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.38%  '
$ws.Range('E7').Value = '  +0.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2661'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06318'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.52%  '
$ws.Range('E10').Value = '  -1.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07732'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.71%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.434'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.61%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.645.00'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.879.28'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5458'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.93%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8229'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.93%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.80'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.65%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.164.19'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.39%  '
$ws.Range('E19').Value = '  +0.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '192.52'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.22%  '
$ws.Range('E22').Value = '  -2.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.092'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.68%  '
$ws.Range('E24').Value = '  +0.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '137.19'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.39%  '
$ws.Range('E26').Value = '  -2.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.226'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.11'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.428'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.06028'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.72%  '
$ws.Range('E31').Value = '  +0.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.558'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.324'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.91%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.648'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9784'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.52%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.411'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.39%  '
$ws.Range('E37').Value = '  -0.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5930'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.96%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01591'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.948'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8628'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.43%  '
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.040.16'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.74%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.792.70'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₈111'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.84%  '
$ws.Range('E47').Value = '  -0.23%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.004'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.62%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.127'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('E50').Value = '  -0.50%  '
